$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column X (col 24) for rows 1-22 with values 0..21
for ($r = 1; $r -le 22; $r++) {
    $ws.Cells.Item($r, 24).Value = $r - 1
}

# Add new row 23 with values 0..22 across columns A-W (col 1-23)
for ($c = 1; $c -le 23; $c++) {
    $ws.Cells.Item(23, $c).Value = $c - 1
}

# Update selection to K24
$ws.Range("K24").Select()

# Update workbook view window position (best effort; xWindow/yWindow track
# the host window placement)
$excel.ActiveWindow.Left = 18820
$excel.ActiveWindow.Top = 1840
